$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.187.54"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "1.685.73"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "215.72"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "0.520"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "23.16"
$ws.Range("E8").Value = "  +8.89%  "
$ws.Range("D9").Value = "0.260"
$ws.Range("E9").Value = "  +3.68%  "
$ws.Range("D10").Value = "0.0626"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "1.921.92"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "1.681.88"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("D14").Value = "4.19"
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("E15").Value = "  +4.05%  "
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "27.170.35"
$ws.Range("D18").Value = "236.20"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("E19").Value = "  -1.80%  "
$ws.Range("D20").Value = "0.0₃0743"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("E23").Value = "  +3.97%  "
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("D25").Value = "146.92"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").Value = "7.33"
$ws.Range("E26").Value = "  +1.28%  "
$ws.Range("D27").Value = "16.43"
$ws.Range("E27").Value = "  +2.33%  "
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("D33").Value = "1.538.71"
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("D34").Value = "3.25"
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("D36").Value = "0.606"
$ws.Range("E36").Value = "  +3.08%  "
$ws.Range("D37").Value = "0.947"
$ws.Range("E37").Value = "  +3.44%  "
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("E40").Value = "  +2.30%  "
$ws.Range("D41").Value = "69.12"
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("D45").Value = "1.830.81"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("E46").Value = "  +1.48%  "
$ws.Range("D47").Value = "90.14"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0112"
$ws.Range("E48").Value = "  +5.05%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.62"
$ws.Range("E49").Value = "  +6.07%  "
$ws.Range("D50").Value = "8.32"
$ws.Range("E50").Value = "  +6.44%  "
$ws.Range("E51").Value = "  -0.31%  "
